$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.321.24"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "2.652.20"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'581.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'144.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "'6.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("D11").Value = "'0.382"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "3.123.85"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "'26.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.09%  "
$ws.Range("D15").Value = "61.234.58"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("E16").Value = "  +3.58%  "
$ws.Range("D17").Value = "2.662.61"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("D18").Value = "'11.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "'4.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").Value = "'355.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").Value = "'6.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'0.524"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("D24").Value = "'64.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'8.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.70%  "
$ws.Range("E28").Value = "  +8.47%  "
$ws.Range("D29").Value = "0.0₃0820"
$ws.Range("E29").Value = "  +3.07%  "
$ws.Range("D30").Value = "'6.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.45%  "
$ws.Range("D31").Value = "'169.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("D33").Value = "'20.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("E34").Value = "  +15.25%  "
$ws.Range("D35").Value = "'4.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.18%  "
$ws.Range("E36").Value = "  +9.92%  "
$ws.Range("D37").Value = "'0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +18.79%  "
$ws.Range("E38").Value = "  +4.80%  "
$ws.Range("D39").Value = "'341.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.14%  "
$ws.Range("D40").Value = "'4.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.07%  "
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").Value = "'5.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.53%  "
$ws.Range("D43").Value = "'0.0579"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.65%  "
$ws.Range("D44").Value = "'20.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.57%  "
$ws.Range("D45").Value = "'21.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("E46").Value = "  +5.31%  "
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "'0.628"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "'0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "2.099.25"
$ws.Range("E51").Value = "  +3.46%  "
